# Auto-generated Word COM-interop script implementing the
# 'Profiles' -> 'Profiles RNS' edit across ORNG_InstallationGuide.docx

function SplitAt($doc, $pos) {
    $r = $doc.Range($pos, $pos)
    $doc.Bookmarks.Add("zzSplitTmp", $r) | Out-Null
    $zb = $doc.Bookmarks.Item("zzSplitTmp")
    $zb.Delete()
}

$d = $word.ActiveDocument

# Remove the pre-existing _GoBack bookmark (currently sits near the
# "OpenSocial support in Profiles." paragraph); it is re-created later
# around the word "Profiles" in the document title, matching the edit.
$hasOldGoBack = $true
try {
    $oldGoBack = $d.Bookmarks.Item("_GoBack")
} catch {
    $hasOldGoBack = $false
}
if ($hasOldGoBack) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# All edits occur in strictly increasing document-order, so scope each
# Find to start where the previous edit ended -- this keeps ambiguous /
# repeated search strings (e.g. "on your Profiles web server" occurs
# both standalone and as a substring of the Tomcat paragraph) pinned to
# the correct, intended occurrence instead of matching earlier text.
$cursor = 0

# --- Edit: title ---
$old_title = "How to Install the Profiles ORNG Extension"
$new_title = "How to Install the Profiles RNS ORNG Extension"
$docEnd_title = $d.Content.End
$rng_title = $d.Range($cursor, $docEnd_title)
$ok_title = $rng_title.Find.Execute($old_title, $false, $false, $false, $false, $false, $true, 1, $false, $new_title, 1)
if (-not $ok_title) { throw "Find failed for title" }
$base_title = $rng_title.Start

$p_title_0 = $base_title + 19
SplitAt $d $p_title_0
$p_title_1 = $base_title + 27
SplitAt $d $p_title_1
$p_title_2 = $base_title + 31
SplitAt $d $p_title_2
$cursor = $rng_title.End

# --- Edit: admin ---
$old_admin = "You need admin rights to the machines running Profiles that you want to apply this to, and you need some amount of time and patience, but that’s pretty much it.  Presumably you know your way around .NET and "
$new_admin = "You need admin rights to the machines running Profiles RNS that you want to apply this to, and you need some amount of time and patience, but that’s pretty much it.  Presumably you know your way around .NET and "
$docEnd_admin = $d.Content.End
$rng_admin = $d.Range($cursor, $docEnd_admin)
$ok_admin = $rng_admin.Find.Execute($old_admin, $false, $false, $false, $false, $false, $true, 1, $false, $new_admin, 1)
if (-not $ok_admin) { throw "Find failed for admin" }
$base_admin = $rng_admin.Start

$p_admin_0 = $base_admin + 55
SplitAt $d $p_admin_0
$p_admin_1 = $base_admin + 59
SplitAt $d $p_admin_1
$cursor = $rng_admin.End

# --- Edit: java ---
$old_java = "on your Profiles web server"
$new_java = "on your Profiles RNS web server"
$docEnd_java = $d.Content.End
$rng_java = $d.Range($cursor, $docEnd_java)
$ok_java = $rng_java.Find.Execute($old_java, $false, $false, $false, $false, $false, $true, 1, $false, $new_java, 1)
if (-not $ok_java) { throw "Find failed for java" }
$base_java = $rng_java.Start

$p_java_0 = $base_java + 17
SplitAt $d $p_java_0
$p_java_1 = $base_java + 21
SplitAt $d $p_java_1
$cursor = $rng_java.End

# --- Edit: tomcat ---
$old_tomcat = "Tomcat on your Profiles web server "
$new_tomcat = "Tomcat on your Profiles RNS web server "
$docEnd_tomcat = $d.Content.End
$rng_tomcat = $d.Range($cursor, $docEnd_tomcat)
$ok_tomcat = $rng_tomcat.Find.Execute($old_tomcat, $false, $false, $false, $false, $false, $true, 1, $false, $new_tomcat, 1)
if (-not $ok_tomcat) { throw "Find failed for tomcat" }
$base_tomcat = $rng_tomcat.Start

$p_tomcat_0 = $base_tomcat + 23
SplitAt $d $p_tomcat_0
$p_tomcat_1 = $base_tomcat + 27
SplitAt $d $p_tomcat_1
$cursor = $rng_tomcat.End

# --- Edit: shindig_host ---
$old_shindig_host = "=<your profiles host machine"
$new_shindig_host = "=<your profiles rns host machine"
$docEnd_shindig_host = $d.Content.End
$rng_shindig_host = $d.Range($cursor, $docEnd_shindig_host)
$ok_shindig_host = $rng_shindig_host.Find.Execute($old_shindig_host, $false, $false, $false, $false, $false, $true, 1, $false, $new_shindig_host, 1)
if (-not $ok_shindig_host) { throw "Find failed for shindig_host" }
$base_shindig_host = $rng_shindig_host.Start

$p_shindig_host_0 = $base_shindig_host + 16
SplitAt $d $p_shindig_host_0
$p_shindig_host_1 = $base_shindig_host + 19
SplitAt $d $p_shindig_host_1
$p_shindig_host_2 = $base_shindig_host + 20
SplitAt $d $p_shindig_host_2
$cursor = $rng_shindig_host.End

# --- Edit: copy_windows ---
$old_copy_windows = "You can then copy this over to your Profiles windows machine, and place it in a /shindig/"
$new_copy_windows = "You can then copy this over to your Profiles RNS windows machine, and place it in a /shindig/"
$docEnd_copy_windows = $d.Content.End
$rng_copy_windows = $d.Range($cursor, $docEnd_copy_windows)
$ok_copy_windows = $rng_copy_windows.Find.Execute($old_copy_windows, $false, $false, $false, $false, $false, $true, 1, $false, $new_copy_windows, 1)
if (-not $ok_copy_windows) { throw "Find failed for copy_windows" }
$base_copy_windows = $rng_copy_windows.Start

$p_copy_windows_0 = $base_copy_windows + 44
SplitAt $d $p_copy_windows_0
$p_copy_windows_1 = $base_copy_windows + 48
SplitAt $d $p_copy_windows_1
$cursor = $rng_copy_windows.End

# --- Edit: orng_db ---
$old_orng_db = " in your Profiles DB"
$new_orng_db = " in your Profiles RNS DB"
$docEnd_orng_db = $d.Content.End
$rng_orng_db = $d.Range($cursor, $docEnd_orng_db)
$ok_orng_db = $rng_orng_db.Find.Execute($old_orng_db, $false, $false, $false, $false, $false, $true, 1, $false, $new_orng_db, 1)
if (-not $ok_orng_db) { throw "Find failed for orng_db" }
$base_orng_db = $rng_orng_db.Start

$p_orng_db_0 = $base_orng_db + 18
SplitAt $d $p_orng_db_0
$p_orng_db_1 = $base_orng_db + 22
SplitAt $d $p_orng_db_1
$cursor = $rng_orng_db.End

# --- Edit: opensocial ---
$old_opensocial = " support in Profiles."
$new_opensocial = " support in Profiles RNS."
$docEnd_opensocial = $d.Content.End
$rng_opensocial = $d.Range($cursor, $docEnd_opensocial)
$ok_opensocial = $rng_opensocial.Find.Execute($old_opensocial, $false, $false, $false, $false, $false, $true, 1, $false, $new_opensocial, 1)
if (-not $ok_opensocial) { throw "Find failed for opensocial" }
$base_opensocial = $rng_opensocial.Start

$p_opensocial_0 = $base_opensocial + 20
SplitAt $d $p_opensocial_0
$p_opensocial_1 = $base_opensocial + 24
SplitAt $d $p_opensocial_1
$cursor = $rng_opensocial.End

# --- Edit: sandbox ---
$old_sandbox = " allows gadgets that are in development to be tested in your profiles environment. For production, you want to leave this commented out. For your development servers you should set it to something private but not to private, as you will want to share this to anyone who wants to work or test gadget changes"
$new_sandbox = " allows gadgets that are in development to be tested in your Profiles RNS environment. For production, you want to leave this commented out. For your development servers you should set it to something private but not to private, as you will want to share this to anyone who wants to work or test gadget changes"
$docEnd_sandbox = $d.Content.End
$rng_sandbox = $d.Range($cursor, $docEnd_sandbox)
$ok_sandbox = $rng_sandbox.Find.Execute($old_sandbox, $false, $false, $false, $false, $false, $true, 1, $false, $new_sandbox, 1)
if (-not $ok_sandbox) { throw "Find failed for sandbox" }
$base_sandbox = $rng_sandbox.Start

$p_sandbox_0 = $base_sandbox + 30
SplitAt $d $p_sandbox_0
$p_sandbox_1 = $base_sandbox + 62
SplitAt $d $p_sandbox_1
$p_sandbox_2 = $base_sandbox + 70
SplitAt $d $p_sandbox_2
$p_sandbox_3 = $base_sandbox + 74
SplitAt $d $p_sandbox_3
$cursor = $rng_sandbox.End

# --- Edit: ucsf ---
$old_ucsf = "You will notice that the default gadgets are hosted on external web sites.  Feel free to use these where they are, or to copy them to one of your own web servers where you can modify them if desired. At UCSF we host them on the same IIS web servers that we use for Profiles."
$new_ucsf = "You will notice that the default gadgets are hosted on external web sites.  Feel free to use these where they are, or to copy them to one of your own web servers where you can modify them if desired. At UCSF we host them on the same IIS web servers that we use for Profiles RNS."
$docEnd_ucsf = $d.Content.End
$rng_ucsf = $d.Range($cursor, $docEnd_ucsf)
$ok_ucsf = $rng_ucsf.Find.Execute($old_ucsf, $false, $false, $false, $false, $false, $true, 1, $false, $new_ucsf, 1)
if (-not $ok_ucsf) { throw "Find failed for ucsf" }
$base_ucsf = $rng_ucsf.Start

$p_ucsf_0 = $base_ucsf + 273
SplitAt $d $p_ucsf_0
$p_ucsf_1 = $base_ucsf + 277
SplitAt $d $p_ucsf_1
$cursor = $rng_ucsf.End

# --- Edit: https ---
$old_https = "You can set the gadgets to run on either HTTPS or HTTP, but not both. For most profiles installations, we expect Profiles to run primarily in HTTP, and gadgets will only show up correctly if a page is rendered as HTTP.  You can, however, have gadgets work in HTTPS if you want to run Profiles in HTTPS only mode.  To do this, set the following:"
$new_https = "You can set the gadgets to run on either HTTPS or HTTP, but not both. For most installations, we expect Profiles RNS to run primarily in HTTP, and gadgets will only show up correctly if a page is rendered as HTTP.  You can, however, have gadgets work in HTTPS if you want to run Profiles RNS in HTTPS only mode.  To do this, set the following:"
$docEnd_https = $d.Content.End
$rng_https = $d.Range($cursor, $docEnd_https)
$ok_https = $rng_https.Find.Execute($old_https, $false, $false, $false, $false, $false, $true, 1, $false, $new_https, 1)
if (-not $ok_https) { throw "Find failed for https" }
$base_https = $rng_https.Start

$p_https_0 = $base_https + 112
SplitAt $d $p_https_0
$p_https_1 = $base_https + 116
SplitAt $d $p_https_1
$p_https_2 = $base_https + 287
SplitAt $d $p_https_2
$p_https_3 = $base_https + 291
SplitAt $d $p_https_3
$cursor = $rng_https.End

# --- Edit: systemdomain ---
$old_systemdomain = " to https for Profiles"
$new_systemdomain = " to https for Profiles RNS"
$docEnd_systemdomain = $d.Content.End
$rng_systemdomain = $d.Range($cursor, $docEnd_systemdomain)
$ok_systemdomain = $rng_systemdomain.Find.Execute($old_systemdomain, $false, $false, $false, $false, $false, $true, 1, $false, $new_systemdomain, 1)
if (-not $ok_systemdomain) { throw "Find failed for systemdomain" }
$base_systemdomain = $rng_systemdomain.Start

$p_systemdomain_0 = $base_systemdomain + 22
SplitAt $d $p_systemdomain_0
$cursor = $rng_systemdomain.End

# Re-create _GoBack bookmark around "Profiles" in the title.
$gbStart = $base_title + 19
$gbEnd = $base_title + 27
$gbRange = $d.Range($gbStart, $gbEnd)
$d.Bookmarks.Add("_GoBack", $gbRange) | Out-Null

Write-Output "All edits applied."
